# Reposition (and, on one slide, resize) the "Oval 8" shape (Id = 9) that
# appears on slides 10-14, matching the author's manual drag/resize edits
# captured in the revisionInfo/changesInfo diff.
#
# PowerPoint's COM object model reports/accepts Shape.Top/Left/Width/Height
# in points, while the underlying OOXML stores offsets/extents in EMU
# (1 pt = 12700 EMU). The host also appears to round-trip these properties
# through a single-precision (f32) float and truncate when converting back
# to EMU, so naive "emu / 12700" point values can land 1 EMU short of the
# intended target. The literals below were chosen so that, after that
# f32 round-trip + truncation, they reproduce the exact target EMU values
# from the diff:
#   slide 10: off (2732313, 3762702) ext (336708, 1145629)
#   slide 11: off (2732313, 3735676) ext (336708, 1219200)  [ext unchanged]
#   slide 12: off (2329215, 3690404) ext (336708, 1219200)  [ext unchanged]
#   slide 13: off (3547078, 3798135) ext (336708, 1219200)  [ext unchanged]
#   slide 14: off (3909599, 3703369) ext (273269, 1219200)  [ext unchanged]
# (Left/Width are untouched in the diff, so they are left as-is.)

$p = $ppt.ActivePresentation

# slideIndex -> @{ Top = <points>; Height = <points or $null> }
$changes = @{
    10 = @{ Top = 296.2757873535156;  Height = 90.20704650878906 }
    11 = @{ Top = 294.14776611328125; Height = $null }
    12 = @{ Top = 290.5830383300781;  Height = $null }
    13 = @{ Top = 299.0657958984375;  Height = $null }
    14 = @{ Top = 291.6039123535156;  Height = $null }
}

foreach ($slideIdx in $changes.Keys) {
    $s = $p.Slides.Item($slideIdx)
    $cfg = $changes[$slideIdx]
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.Id -eq 9) {
            $sh.Top = $cfg.Top
            if ($null -ne $cfg.Height) {
                $sh.Height = $cfg.Height
            }
        }
    }
}
